# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G (header "K") values for rows 2-15 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 3
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 2
    12 = 1
    13 = 2
    14 = 2
    15 = 2
}

foreach ($row in $newKValues.Keys) {
    $ws.Range("G$row").Value = $newKValues[$row]
}
